{"js": "// Change \"Database: MySQL for data storage.\" to\n// \"Database: MongoDB for data storage.\" by replacing the \"MySQL\" token\n// with \"MongoDB\" in the document body.\n\nconst results = context.document.body.search(\"MySQL\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find 'MySQL' text to replace.\");\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"MongoDB\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the \"Database\" bullet so it references MongoDB instead of MySQL:\n#   \"Database: MySQL for data storage.\" -> \"Database: MongoDB for data storage.\"\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"MySQL\", $false, $false, $false, $false, $false, $true, 1, $false, \"MongoDB\", 2)\n"}
